# Auto-generated script to apply scheduled-runner price/profit updates
# across the FFXIV crafting-leve profit tables (ALC, ARM, BSM, CRP, CUL, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("N88").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("N91").Value = 0
$ws.Range("H100").Value = 2907.25
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H113").Value = 7482.7144
$ws.Range("I113").Value = 7666.5
$ws.Range("J113").Value = 6380
$ws.Range("K113").Value = 7666.5
$ws.Range("L113").Value = 6380
$ws.Range("M113").Value = -4412.5
$ws.Range("N113").Value = -12888
$ws.Range("H129").Value = 1056
$ws.Range("I129").Value = 1056
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 3168
$ws.Range("L129").Value = 0
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = 1832
$ws.Range("H137").Value = 4532.8887
$ws.Range("I137").Value = 4179.4
$ws.Range("J137").Value = 4974.75
$ws.Range("K137").Value = 12538.2
$ws.Range("L137").Value = 14924.25
$ws.Range("M137").Value = -9988.199999999999
$ws.Range("N137").Value = -20024.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5477.25
$ws.Range("I2").Value = 5477.25
$ws.Range("K2").Value = 5477.25
$ws.Range("M2").Value = -5364.25
$ws.Range("H15").Value = 50000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 50000
$ws.Range("K15").Value = 0
$ws.Range("L15").ClearContents()
$ws.Range("M15").Value = 50000
$ws.Range("N15").Value = -50700
$ws.Range("H45").Value = 3651.875
$ws.Range("I45").Value = 3916.4285
$ws.Range("J45").Value = 1800
$ws.Range("K45").Value = 3916.4285
$ws.Range("L45").Value = 1800
$ws.Range("M45").Value = -3539.4285
$ws.Range("N45").Value = -2554
$ws.Range("H74").Value = 20814.889
$ws.Range("I74").Value = 20280.12
$ws.Range("K74").Value = 20280.12
$ws.Range("M74").Value = -19406.12
$ws.Range("H77").Value = 20814.889
$ws.Range("I77").Value = 20280.12
$ws.Range("K77").Value = 101400.6
$ws.Range("M77").Value = -97032.59999999999
$ws.Range("H110").Value = 5550
$ws.Range("I110").Value = 2771.4285
$ws.Range("J110").Value = 25000
$ws.Range("K110").Value = 2771.4285
$ws.Range("L110").Value = 25000
$ws.Range("M110").Value = -726.4285
$ws.Range("N110").Value = -29090
$ws.Range("H116").Value = 5477.25
$ws.Range("I116").Value = 5477.25
$ws.Range("K116").Value = 5477.25
$ws.Range("M116").Value = -3183.25
$ws.Range("H122").Value = 3989
$ws.Range("I122").Value = 4662.3335
$ws.Range("J122").Value = 1969
$ws.Range("K122").Value = 13987.0005
$ws.Range("L122").Value = 5907
$ws.Range("M122").Value = -11537.0005
$ws.Range("N122").Value = -10807

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5477.25
$ws.Range("I3").Value = 5477.25
$ws.Range("K3").Value = 5477.25
$ws.Range("M3").Value = -5363.25
$ws.Range("H8").Value = 500
$ws.Range("J8").Value = 500
$ws.Range("L8").Value = 500
$ws.Range("N8").Value = -780
$ws.Range("H94").Value = 3234.8333
$ws.Range("I94").Value = 5454.5
$ws.Range("K94").Value = 5454.5
$ws.Range("M94").Value = -5003.5
$ws.Range("H99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("N99").Value = 0
$ws.Range("H107").Value = 999.6667
$ws.Range("J107").Value = 1200
$ws.Range("L107").Value = 1200
$ws.Range("N107").Value = -5040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2337.25
$ws.Range("I31").Value = 2233.3333
$ws.Range("J31").Value = 2399.6
$ws.Range("K31").Value = 2233.3333
$ws.Range("L31").Value = 2399.6
$ws.Range("M31").Value = -1938.3333
$ws.Range("N31").Value = -2989.6
$ws.Range("H34").Value = 2337.25
$ws.Range("I34").Value = 2233.3333
$ws.Range("J34").Value = 2399.6
$ws.Range("K34").Value = 2233.3333
$ws.Range("L34").Value = 2399.6
$ws.Range("M34").Value = -2031.3333
$ws.Range("N34").Value = -2803.6
$ws.Range("H58").Value = 4715.273
$ws.Range("I58").Value = 4715.273
$ws.Range("K58").Value = 4715.273
$ws.Range("M58").Value = -4512.273
$ws.Range("H110").Value = 47000
$ws.Range("J110").Value = 47000
$ws.Range("L110").Value = 47000
$ws.Range("N110").Value = -55180
$ws.Range("H134").Value = 4187.3887
$ws.Range("I134").Value = 3822.1538
$ws.Range("J134").Value = 5137
$ws.Range("K134").Value = 11466.4614
$ws.Range("L134").Value = 15411
$ws.Range("M134").Value = -8931.4614
$ws.Range("N134").Value = -20481
$ws.Range("H136").Value = 4715.273
$ws.Range("I136").Value = 4715.273
$ws.Range("K136").Value = 14145.819
$ws.Range("M136").Value = -11595.819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 135.42857
$ws.Range("I2").Value = 106
$ws.Range("J2").Value = 174.66667
$ws.Range("K2").Value = 636
$ws.Range("L2").Value = 1048.00002
$ws.Range("M2").Value = -523
$ws.Range("N2").Value = -1274.00002
$ws.Range("H4").Value = 4347.857
$ws.Range("I4").Value = 5056.6665
$ws.Range("K4").Value = 15169.9995
$ws.Range("M4").Value = -15057.9995
$ws.Range("H18").Value = 1511.5385
$ws.Range("I18").Value = 2172.5
$ws.Range("K18").Value = 6517.5
$ws.Range("M18").Value = -6348.5
$ws.Range("H38").Value = 301
$ws.Range("I38").Value = 42.5
$ws.Range("K38").Value = 127.5
$ws.Range("M38").Value = 219.5
$ws.Range("H99").Value = 4949
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 4949
$ws.Range("K99").Value = 0
$ws.Range("L99").ClearContents()
$ws.Range("M99").Value = 14847
$ws.Range("N99").Value = -19339
$ws.Range("H122").Value = 3234.75
$ws.Range("J122").Value = 3699.2
$ws.Range("L122").Value = 33292.8
$ws.Range("N122").Value = -38192.8
$ws.Range("H123").Value = 4952.381
$ws.Range("I123").Value = 4000
$ws.Range("K123").Value = 12000
$ws.Range("M123").Value = -9550
$ws.Range("H126").Value = 6909
$ws.Range("I126").Value = 999
$ws.Range("K126").Value = 2997
$ws.Range("M126").Value = 1943
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").ClearContents()
$ws.Range("N127").Value = 0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1989.9166
$ws.Range("I55").Value = 2137.9
$ws.Range("J55").Value = 1250
$ws.Range("K55").Value = 2137.9
$ws.Range("L55").Value = 1250
$ws.Range("M55").Value = -1964.9
$ws.Range("N55").Value = -1596

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 11501
$ws.Range("J82").Value = 11501
$ws.Range("L82").Value = 11501
$ws.Range("N82").Value = -12267
$ws.Range("H85").Value = 11501
$ws.Range("J85").Value = 11501
$ws.Range("L85").Value = 11501
$ws.Range("N85").Value = -14153
$ws.Range("H113").Value = 388.36365
$ws.Range("I113").Value = 391.33334
$ws.Range("K113").Value = 1174.00002
$ws.Range("M113").Value = 995.9999800000001
$ws.Range("H136").Value = 1624.9642
$ws.Range("I136").Value = 1624.9642
$ws.Range("K136").Value = 4874.892599999999
$ws.Range("M136").Value = -2324.892599999999
